# "Generate Report for Archive"
#
# 1) The localization status for the sample file moved from
#    "Ready for handoff" to "In Translation". That text lives in cell
#    E2/F2 on the "Overview" sheet and in C2 on each language sheet
#    ("zh-cn", "de-de") - all four cells share the same string, so
#    rewriting it anywhere it appears covers every occurrence.
# 2) The "Status" column got narrower: column E & F on "Overview" and
#    column C on each language sheet shrank from ~17.22 characters to
#    ~13.41 characters.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# ColumnWidth is expressed in characters and Excel snaps whatever is
# assigned to the nearest whole pixel, so we pick the input that lands on
# the stored width closest to the narrower target width used in the
# updated report.
$newStatusColumnWidth = 12.5

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    for ($r = 1; $r -le $used.Rows.Count; $r++) {
        for ($c = 1; $c -le $used.Columns.Count; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            # Compare with the literal on the left so PowerShell's -eq
            # doesn't coerce the string into the cell's native type (e.g.
            # a boolean "True" cell would otherwise equal any non-empty
            # string if the cell's value were the left-hand operand).
            if ($oldStatus -eq $cell.Value2) {
                $cell.Value2 = $newStatus
            }
        }
    }
}

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E1").ColumnWidth = $newStatusColumnWidth
$overview.Range("F1").ColumnWidth = $newStatusColumnWidth

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C1").ColumnWidth = $newStatusColumnWidth

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C1").ColumnWidth = $newStatusColumnWidth
